# Commit: "Add functionality for adding groups to the db"
# Populate the (previously empty) "Groups" worksheet with the new group rows.
$wb = $excel.ActiveWorkbook
$groups = $wb.Worksheets.Item("Groups")

$groups.Range("A1").Value = "okokokok"
$groups.Range("B1").Value = -615761128

$groups.Range("A2").Value = "okokokok"
$groups.Range("B2").Value = -615761128

$groups.Range("A3").Value = "okokokok"
$groups.Range("B3").Value = -615761128

$groups.Range("A4").Value = "annie's test"
$groups.Range("B4").Value = -485430438
